$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New game rows to append (Away team, Away Pts, Home team, Home Pts, Overtime,
# Attend., Arena, Win, Loss) matching the sheet's existing header layout.
$rows = @(
    @("Dallas Mavericks", 119, "Cleveland Cavaliers", 121, "No", 17832, "Rocket Mortgage Fieldhouse", "Cleveland Cavaliers", "Dallas Mavericks"),
    @("Brooklyn Nets", 81, "Orlando Magic", 108, "No", 17832, "Amway Center", "Orlando Magic", "Brooklyn Nets"),
    @("Golden State Warriors", 123, "Washington Wizards", 112, "No", 17832, "Capital One Arena", "Golden State Warriors", "Washington Wizards"),
    @("Utah Jazz", 97, "Atlanta Hawks", 124, "No", 17832, "State Farm Arena", "Atlanta Hawks", "Utah Jazz"),
    @("Philadelphia 76ers", 99, "Boston Celtics", 117, "No", 17832, "TD Garden", "Boston Celtics", "Philadelphia 76ers"),
    @("San Antonio Spurs", 105, "Minnesota Timberwolves", 114, "No", 17832, "Target Center", "Minnesota Timberwolves", "San Antonio Spurs"),
    @("New Orleans Pelicans", 115, "New York Knicks", 92, "No", 17832, "Madison Square Garden (IV)", "New Orleans Pelicans", "New York Knicks"),
    @("Detroit Pistons", 105, "Chicago Bulls", 95, "No", 17832, "United Center", "Detroit Pistons", "Chicago Bulls"),
    @("Charlotte Hornets", 85, "Milwaukee Bucks", 123, "No", 17832, "Fiserv Forum", "Milwaukee Bucks", "Charlotte Hornets"),
    @("Houston Rockets", 95, "Oklahoma City Thunder", 112, "No", 17832, "Paycom Center", "Oklahoma City Thunder", "Houston Rockets"),
    @("Miami Heat", 106, "Portland Trail Blazers", 96, "No", 17832, "Moda Center", "Miami Heat", "Portland Trail Blazers")
)

$startRow = $ws.UsedRange.Rows.Count + 1
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $targetRow = $startRow + $i
    for ($c = 0; $c -lt $r.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $r[$c]
    }
}

$lastRow = $startRow + $rows.Length - 1
$lastCell = $ws.Cells.Item($lastRow, 1)
[void]$lastCell.Select()
$excel.ActiveWindow.ScrollRow = $lastRow - 32
